$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = 0 (numeric)
$ws.Range("B1").Value = 0

# A2 = 0 (numeric)
$ws.Range("A2").Value = 0

# B2 = "disconnected_elements" (text -> shared string), no special style
$ws.Range("B2").Value = "disconnected_elements"

# Style B1: bold font, thin border all around, centered horizontally, top vertically
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Borders.LineStyle = 1        # xlContinuous
$r1.Borders.Weight = 2           # xlThin
$r1.Borders.ColorIndex = -4105   # xlColorIndexAutomatic

# Copy B1's format onto A2 so both cells share the exact same style record
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
